$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1,1).Value = "Instance Name"
$ws.Cells.Item(1,2).Value = "Instance Type"
$ws.Cells.Item(1,3).Value = "Region"
$ws.Cells.Item(1,4).Value = "OS"
$ws.Cells.Item(1,5).Value = "Storage"
$ws.Cells.Item(1,6).Value = "Storage Type"

# Data rows - one row per instance (Instance Name replaces Count column)
$rows = @(
    @("web-server-01",   "t2.micro",    "us-east-1",      "Linux",   100, "SSD"),
    @("web-server-02",   "t2.micro",    "us-east-1",      "Linux",   100, "SSD"),
    @("app-server-01",   "m5.xlarge",   "us-west-2",      "Windows", 500, "SSD"),
    @("api-server-01",   "c5.2xlarge",  "ap-southeast-1", "Linux",   200, "SSD"),
    @("api-server-02",   "c5.2xlarge",  "ap-southeast-1", "Linux",   200, "SSD"),
    @("api-server-03",   "c5.2xlarge",  "ap-southeast-1", "Linux",   200, "SSD"),
    @("db-server-01",    "r5.large",    "eu-west-1",      "Linux",   300, "SSD"),
    @("db-server-02",    "r5.large",    "eu-west-1",      "Linux",   300, "SSD"),
    @("cache-server-01", "m6i.2xlarge", "us-east-1",      "Linux",   250, "SSD")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
    $r++
}
